# xp_SkillQuest_2026-02-18.xlsx — 2026-02-18 roster upload
#
# The "Date" column (D2:D216) was entered/re-entered as a literal text
# string ("2026-02-18") instead of a real date serial, so the whole
# column is retyped as Text before the value is (re)written — otherwise
# Excel would just re-parse the string back into the existing date serial
# using the column's current date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = $ws.Range("D2:D216")

# Switch the column to Text format *before* writing the string so Excel
# stores it verbatim instead of reinterpreting it as a date serial.
$dateCol.NumberFormat = "@"
$dateCol.Value = "2026-02-18"

# Reflect the in-progress selection/scroll left behind after editing the
# column: the last edited cell (D215) is active within the D2:D215
# selection, scrolled down so row 196 is at the top of the window.
$ws.Range("D2:D215").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 196
$win.ScrollColumn = 1
